$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44496
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14520
$ws.Range("P2").Value = 581
# Row 3
$ws.Range("D3").Value = 44512
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("N3").Value = "`$/saco 25 kilos"
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 580
# Row 4
$ws.Range("D4").Value = 44657
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24400
$ws.Range("O4").Value = "Carahue"
$ws.Range("P4").Value = 976
# Row 5
$ws.Range("D5").Value = 44482
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 24385
$ws.Range("O5").Value = "Región de O'Higgins"
$ws.Range("P5").Value = 975
# Row 6
$ws.Range("D6").Value = 44540
$ws.Range("J6").Value = 110
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 16545
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 662
# Row 7
$ws.Range("D7").Value = 44615
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Carahue"
$ws.Range("P7").Value = 1160
# Row 8
$ws.Range("D8").Value = 44659
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 24571
$ws.Range("O8").Value = "Carahue"
$ws.Range("P8").Value = 983
# Row 9
$ws.Range("D9").Value = 44532
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14400
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 576
# Row 10
$ws.Range("D10").Value = 44595
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 26000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 27200
$ws.Range("N10").Value = "`$/saco 25 kilos"
$ws.Range("O10").Value = "Carahue"
$ws.Range("P10").Value = 1088
# Row 11
$ws.Range("D11").Value = 44545
$ws.Range("J11").Value = 180
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15444
$ws.Range("O11").Value = "Carahue"
$ws.Range("P11").Value = 618
# Row 12
$ws.Range("D12").Value = 44518
$ws.Range("J12").Value = 350
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14571
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 583
# Row 13
$ws.Range("D13").Value = 44517
$ws.Range("J13").Value = 110
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17455
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 698
# Row 14
$ws.Range("D14").Value = 44539
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13400
$ws.Range("P14").Value = 536
# Row 15
$ws.Range("D15").Value = 44162
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17500
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("P15").Value = 700
# Row 16
$ws.Range("D16").Value = 44643
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 26000
$ws.Range("M16").Value = 25444
$ws.Range("O16").Value = "Carahue"
$ws.Range("P16").Value = 1018
# Row 17
$ws.Range("D17").Value = 44589
$ws.Range("H17").Value = "Perfection"
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 22000
$ws.Range("L17").Value = 23000
$ws.Range("M17").Value = 22500
$ws.Range("N17").Value = "`$/malla 25 kilos"
$ws.Range("P17").Value = 900
# Row 18
$ws.Range("D18").Value = 44483
$ws.Range("J18").Value = 220
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19455
$ws.Range("N18").Value = "`$/saco 25 kilos"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 778
# Row 19
$ws.Range("D19").Value = 44673
$ws.Range("J19").Value = 220
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 26000
$ws.Range("M19").Value = 25455
$ws.Range("O19").Value = "Carahue"
$ws.Range("P19").Value = 1018
# Row 20
$ws.Range("D20").Value = 44503
$ws.Range("H20").Value = "Perfection"
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("N20").Value = "`$/malla 25 kilos"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 620
# Row 21
$ws.Range("D21").Value = 44335
$ws.Range("K21").Value = 30000
$ws.Range("L21").Value = 32000
$ws.Range("M21").Value = 31000
$ws.Range("N21").Value = "`$/malla 25 kilos"
$ws.Range("O21").Value = "Provincia de Huasco"
$ws.Range("P21").Value = 1240
# Row 22
$ws.Range("D22").Value = 44328
$ws.Range("H22").Value = "Perfection"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 33000
$ws.Range("L22").Value = 34000
$ws.Range("M22").Value = 33500
$ws.Range("N22").Value = "`$/malla 25 kilos"
$ws.Range("O22").Value = "Provincia de Huasco"
$ws.Range("P22").Value = 1340
# Row 23
$ws.Range("D23").Value = 44399
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 39000
$ws.Range("L23").Value = 40000
$ws.Range("M23").Value = 39600
$ws.Range("N23").Value = "`$/malla 25 kilos"
$ws.Range("O23").Value = "Provincia de Huasco"
$ws.Range("P23").Value = 1584
# Row 24
$ws.Range("D24").Value = 44454
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 36000
$ws.Range("L24").Value = 38000
$ws.Range("M24").Value = 37000
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 1480
# Row 25
$ws.Range("D25").Value = 44342
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 30000
$ws.Range("L25").Value = 32000
$ws.Range("M25").Value = 31000
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 1240
# Row 26
$ws.Range("D26").Value = 44519
$ws.Range("H26").Value = "Perfection"
$ws.Range("J26").Value = 240
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17583
$ws.Range("P26").Value = 703
# Row 27
$ws.Range("D27").Value = 44505
$ws.Range("J27").Value = 210
$ws.Range("K27").Value = 6500
$ws.Range("L27").Value = 7000
$ws.Range("M27").Value = 6714
$ws.Range("N27").Value = "`$/malla 25 kilos"
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 269
# Row 28
$ws.Range("D28").Value = 44671
$ws.Range("J28").Value = 110
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = 25545
$ws.Range("N28").Value = "`$/saco 25 kilos"
$ws.Range("O28").Value = "Carahue"
$ws.Range("P28").Value = 1022
# Row 29
$ws.Range("D29").Value = 44631
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 25000
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 24467
$ws.Range("P29").Value = 979
# Row 30
$ws.Range("D30").Value = 44533
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 14000
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = 14375
$ws.Range("N30").Value = "`$/malla 25 kilos"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 575
# Row 31
$ws.Range("D31").Value = 44629
$ws.Range("J31").Value = 35
$ws.Range("K31").Value = 25000
$ws.Range("L31").Value = 26000
$ws.Range("M31").Value = 25429
$ws.Range("N31").Value = "`$/saco 25 kilos"
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 1017
